$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.258.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.40%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.785.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.63%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3784"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.76%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3430"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.93%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.39"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.197"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07490"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.472"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.787.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.093"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001097"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.87%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06648"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.649"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.259.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.406"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.494"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.545"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.991.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "133.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.58%  "

$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.024"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.37%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.090"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08701"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.90%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.664"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6940"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.447"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2204"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.832"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06322"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02336"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.239"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.91%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6525"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.842"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.150"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07135"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.48%  "
